$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "35.749.69"
$ws.Range("E2").Value = "  +0.67%  "
$ws.Range("D3").Value = "1.899.84"
$ws.Range("E3").Value = "  +0.26%  "
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").Value = "247.23"
$ws.Range("E5").Value = "  -0.33%  "
$ws.Range("D6").Value = "0.692"
$ws.Range("E6").Value = "  +0.15%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").Value = "43.30"
$ws.Range("E8").Value = "  -1.76%  "
$ws.Range("D9").Value = "57.14"
$ws.Range("E9").Value = "  +9.85%  "
$ws.Range("E10").Value = "  +1.38%  "
$ws.Range("D11").Value = "0.0757"
$ws.Range("E11").Value = "  +2.11%  "
$ws.Range("D12").Value = "0.0986"
$ws.Range("E12").Value = "  +1.74%  "
$ws.Range("D13").Value = "14.59"
$ws.Range("E13").Value = "  +11.38%  "
$ws.Range("D14").Value = "0.804"
$ws.Range("E14").Value = "  +11.11%  "
$ws.Range("D15").Value = "2.176.78"
$ws.Range("E15").Value = "  +0.23%  "
$ws.Range("E16").Value = "  +2.32%  "
$ws.Range("D17").Value = "1.904.27"
$ws.Range("E17").Value = "  +1.11%  "
$ws.Range("D18").Value = "35.706.10"
$ws.Range("E18").Value = "  +0.51%  "
$ws.Range("D19").Value = "73.80"
$ws.Range("E19").Value = "  +0.15%  "
$ws.Range("D20").Value = "0.0₃0833"
$ws.Range("E20").Value = "  +1.27%  "
$ws.Range("D21").Value = "247.13"
$ws.Range("E21").Value = "  -0.38%  "
$ws.Range("D22").Value = "13.05"
$ws.Range("E23").Value = "  +4.84%  "
$ws.Range("D24").Value = "2.69"
$ws.Range("E24").Value = "  +6.24%  "
$ws.Range("E25").Value = "  +0.01%  "
$ws.Range("D26").Value = "2.15"
$ws.Range("E26").Value = "  -1.94%  "
$ws.Range("D27").Value = "167.00"
$ws.Range("E27").Value = "  +0.90%  "
$ws.Range("D28").Value = "8.69"
$ws.Range("E28").Value = "  +2.54%  "
$ws.Range("D29").Value = "18.42"
$ws.Range("E29").Value = "  +0.36%  "
$ws.Range("E30").Value = "  +0.52%  "
$ws.Range("B31").Value = "Hedera"
$ws.Range("C31").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D31").Value = "0.0607"
$ws.Range("E31").Value = "  +4.56%  "
$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").Value = "4.37"
$ws.Range("E32").Value = "  +2.85%  "
$ws.Range("D33").Value = "4.27"
$ws.Range("E33").Value = "  +0.82%  "
$ws.Range("E34").Value = "  +15.78%  "
$ws.Range("E35").Value = "  -0.06%  "
$ws.Range("D36").Value = "1.49"
$ws.Range("E36").Value = "  -16.35%  "
$ws.Range("E37").Value = "  +0.04%  "
$ws.Range("D38").Value = "1.97"
$ws.Range("E38").Value = "  -2.18%  "
$ws.Range("E39").Value = "  +7.06%  "
$ws.Range("D40").Value = "0.0229"
$ws.Range("E40").Value = "  +7.44%  "
$ws.Range("D41").Value = "99.30"
$ws.Range("E41").Value = "  +1.73%  "
$ws.Range("D42").Value = "17.10"
$ws.Range("E42").Value = "  -0.11%  "
$ws.Range("B43").Value = "ARBITRUM"
$ws.Range("C43").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D43").Value = "1.09"
$ws.Range("E43").Value = "  -0.19%  "
$ws.Range("B44").Value = "Gas"
$ws.Range("C44").Value = "https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas"
$ws.Range("D44").Value = "14.00"
$ws.Range("E44").Value = "  +15.39%  "
$ws.Range("D45").Value = "1.320.43"
$ws.Range("E45").Value = "  +2.24%  "
$ws.Range("D46").Value = "2.37"
$ws.Range("E46").Value = "  -0.20%  "
$ws.Range("E47").Value = "  +1.65%  "
$ws.Range("E48").Value = "  -0.17%  "
$ws.Range("E49").Value = "  +0.25%  "
$ws.Range("D50").Value = "6.45"
$ws.Range("E50").Value = "  +0.87%  "
$ws.Range("D51").Value = "42.65"
$ws.Range("E51").Value = "  -1.31%  "
